$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G4").Value = "2016-08-31 02:52:42"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H4").Value = "2016-08-31 02:52:38"
$wsZhCn.Range("K4").Value = "2016-08-31 02:52:56"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H4").Value = "2016-08-31 02:52:42"
$wsDeDe.Range("K4").Value = "2016-08-31 02:53:10"
